$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","I","J","K","L","N","O")

$newValues = @{
    2 = @(12.44920852904582,8.837927998651603,9.6366667936373,32.05724342210324,32.70405752597234,15.41561275117694,24.36903604213749,10.44565769584516,9.307327154608577,11.58239160488619,19.32661972788484,23.96315434154425)
    3 = @(12.21325946675635,8.827412745671372,9.61838139383401,32.1305562419502,32.79762125717165,15.45856872640799,24.4521700932462,10.46641076717269,9.135829970393692,11.57677813002858,19.37686786354418,24.03641297067031)
    4 = @(12.06776280606209,8.8211198334429,9.6085594332074,32.18161577542555,32.86286376587093,15.48687989552123,24.50688009120664,10.48000882624959,9.029937207802414,11.57480983969698,19.40925856471038,24.085349776001)
    5 = @(12.00839281198697,8.818597299771492,9.604913526978283,32.20394201747429,32.89140625383465,15.49890424601764,24.5300970257092,10.48576575242839,8.986689964028047,11.57438124951328,19.42284582793528,24.10628638787159)
    6 = @(11.99853200824877,8.81818098658243,9.604329755996469,32.20774098659368,32.89626369776209,15.50093032759489,24.53400789102527,10.48673472135047,8.979504650439754,11.57433268951672,19.42512543399195,24.10982295259306)
    7 = @(12.06696233713975,8.821085643096357,9.608508815186557,32.18191072580117,32.86324078737648,15.48704008632662,24.5071894684333,10.48008559256843,9.02935427312711,11.57480254506551,19.40944023580197,24.08562810822745)
    8 = @(12.36802864986003,8.834268866714007,9.630072318331855,32.08126644193985,32.73469821090976,15.43002244408705,24.39694020953759,10.45263605507125,9.248349944947682,11.58015032180174,19.3436266538899,23.98759277922796)
    9 = @(12.95016453644861,8.861385332062612,9.683357493789266,31.9319090433305,32.54464473697322,15.33355052841044,24.20980279974063,10.40557555934163,9.670796675446077,11.60228654246161,19.22672421022177,23.8267422323383)
    10 = @(13.36857497139917,8.882035034597202,9.728994814195191,31.85148676167984,32.44304174845323,15.27199368476353,24.08999890392342,10.37509754952846,9.973943952575002,11.62553611898163,19.14818158082594,23.7277176380599)
    11 = @(13.55606572717175,8.891577616757683,9.751114484667587,31.82127017667716,32.40511716108526,15.24600752211296,24.03933119697284,10.36211597265929,10.10970100030499,11.63760243775146,19.11403150448939,23.68682975272323)
    12 = @(13.62658912414226,8.895211710362943,9.759681376696825,31.81074360873189,32.3919514715435,15.236456725384,24.02069527008172,10.35732670354262,10.16075458003227,11.64238330276036,19.10132584155258,23.67194480859526)
    13 = @(13.61142279713853,8.894428144748831,9.757827942113122,31.81296996447889,32.39473372878951,15.23850078588227,24.02468435425112,10.35835253541798,10.14977575762101,11.6413442900662,19.10405218565642,23.675123927703)
    14 = @(13.56187761925795,8.89187618962289,9.751815495149044,31.82038579605581,32.40401003787161,15.24521597031204,24.03778696936991,10.36171942221799,10.11390857161848,11.63799153702741,19.11298167296911,23.68559316439503)
    15 = @(13.53146597799224,8.890315689716424,9.748157381157489,31.82504746914955,32.40984781047356,15.24936691869997,24.04588442987493,10.36379820908291,10.09189137401053,11.63596535619066,19.11848067766482,23.6920838151771)
    16 = @(13.356259007269,8.881414325355214,9.727576185579851,31.85358963987031,32.44568736971372,15.27373248577894,24.09338724725739,10.37596366005574,9.965024728076155,11.62477731595874,19.15044507386783,23.73047347727104)
    17 = @(13.24799805282733,8.875991140809555,9.715295039198423,31.87273046810693,32.46980042341168,15.28919615818356,24.123509957229,10.38365264398784,9.886612880310723,11.6182935973873,19.17045807743787,23.75508975645037)
    18 = @(13.18546502292249,8.87288594046041,9.708359465882655,31.88433915542693,32.48445025712824,15.29828026302157,24.14119640826281,10.38815828246412,9.841313044969784,11.61470474436433,19.18211773643435,23.76963982461823)
    19 = @(13.16424911134907,8.871837016796649,9.706033363307716,31.88837259575468,32.4895444255808,15.30138859553515,24.14724668628137,10.38969810572075,9.825942494999246,11.61351381421136,19.18609106450692,23.77463344541635)
    20 = @(13.25955044931541,8.876566998086144,9.716589152100347,31.87063085816904,32.46715273305416,15.2875303818387,24.12026601704416,10.38282553697045,9.89498093582214,11.61896928704551,19.16831227466627,23.75242879664798)
    21 = @(13.57644364499101,8.892625210824846,9.753576361382715,31.81818273264718,32.40125289548636,15.24323570100144,24.03392346540814,10.3607270540269,10.12445361324591,11.63897060056077,19.11035273323456,23.68250185236046)
    22 = @(13.78075346719264,8.903239393971782,9.778858551757461,31.78924266986216,32.36515292329985,15.21597441197933,23.98070423571667,10.34702199264404,10.27234049075703,11.65327472451973,19.07379106297008,23.64028869247813)
    23 = @(13.67198631319161,8.897563794088038,9.765265136114259,31.80420015281114,32.38378168619806,15.23036995166557,24.0088146168984,10.35426928668922,10.19361595863742,11.64552851664824,19.09318438308006,23.66249937326775)
    24 = @(13.25432851645333,8.876306613513254,9.716003693880539,31.87157820928303,32.46834730290223,15.28828287552137,24.12173145453325,10.38319920682493,9.891198419082365,11.61866337545106,19.16928191259206,23.75363057736143)
    25 = @(12.79401567073515,8.853918704281485,9.667786334014993,31.96716925100297,32.5893953717995,15.35800968955083,24.25732044830302,10.41758508272529,9.557576351679632,11.5950623317521,19.25705472269057,23.86689400385482)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

